$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-171
# from serial date 45204 (2023-10-05) to 45205 (2023-10-06).
$ws.Range("C2:C171").Value = 45205
